$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'66.497.44"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  +3.10%  "
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'3.255.48"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'  +6.44%  "
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Value = "'  -0.04%  "
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'583.52"
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = "'153.71"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'  +7.22%  "
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = "'  -0.02%  "
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'3.245.78"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  +6.36%  "
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'0.518"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  +5.25%  "
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'7.10"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  +8.49%  "
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.Value = "'  +5.58%  "
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'0.492"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  +4.41%  "
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'38.09"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  +2.61%  "
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.Value = "'  +5.27%  "
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'3.776.09"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  +6.54%  "
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'558.20"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  +12.21%  "
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = "'66.545.56"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'  +3.02%  "
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = "'3.254.77"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  +6.49%  "
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.Value = "'  +2.97%  "
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'7.17"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  +5.66%  "
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'14.61"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  +5.20%  "
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.Value = "'  +7.89%  "
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'7.87"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  +9.15%  "
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'13.62"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'  +6.50%  "
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'82.02"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  +3.00%  "
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Value = "'1.00"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  +0.01%  "
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.Value = "'9.32"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  +17.99%  "
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.Value = "'  +8.28%  "
$c.Style = "Normal"

$c = $ws.Range("E29")
$c.Value = "'  +5.49%  "
$c.Style = "Normal"

$c = $ws.Range("D30")
$c.Value = "'27.90"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'  +6.07%  "
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.Value = "'2.78"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  +5.61%  "
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.Value = "'  -0.18%  "
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.Value = "'  +5.14%  "
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.Value = "'568.46"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  +8.00%  "
$c.Style = "Normal"

$c = $ws.Range("D35")
$c.Value = "'5.74"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  +3.52%  "
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.Value = "'6.44"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'  +6.67%  "
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.Value = "'55.39"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  +4.29%  "
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.Value = "'  +11.36%  "
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = "'0.0869"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  +7.20%  "
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'0.132"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  +7.39%  "
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'3.04"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'  +10.40%  "
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'3.216.99"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  +10.72%  "
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'8.69"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  +3.20%  "
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'0.284"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  +14.23%  "
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.Value = "'  +9.40%  "
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.Value = "'26.65"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'  +5.77%  "
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'0.0₃0563"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  +3.33%  "
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.Value = "'  +0.06%  "
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.Value = "'126.49"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'  +4.50%  "
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.Value = "'  +3.52%  "
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.Value = "'  +7.41%  "
$c.Style = "Normal"
